$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows after row 3 (pushes old rows 4-29 down to rows 6-31)
$ws.Rows("4:5").Insert()

# Fix formatting on new A4:A5 cells to match the bold+bordered style used elsewhere in column A
$ws.Range("A4:A5").Font.Bold = $true
$ws.Range("A4:A5").Borders.LineStyle = 1
$ws.Range("A4:A5").HorizontalAlignment = -4108
$ws.Range("A4:A5").VerticalAlignment = -4160

# New row 4: Holden
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.8802992153299319
$ws.Range("D4").Value = 1.103694516580722
$ws.Range("E4").Value = 0.9881753279078053
$ws.Range("F4").Value = 1.007981965640347
$ws.Range("G4").Value = 0.8802992153299319
$ws.Range("H4").Value = 0.9938876833932523
$ws.Range("I4").Value = 0.9224722719409898
$ws.Range("J4").Value = 0.9881753279078053
$ws.Range("K4").Value = 0.9881753279078053
$ws.Range("L4").Value = 1.061797652766367
$ws.Range("M4").Value = 0.9655800664465234
$ws.Range("N4").Value = 0.9881753279078053
$ws.Range("O4").Value = 1.103694516580722
$ws.Range("P4").Value = 0.9919968659553271
$ws.Range("Q4").Value = 1.034637291513623
$ws.Range("R4").Value = 0.9907230199394865
$ws.Range("S4").Value = 0.9831912661190593
$ws.Range("T4").Value = 0.9907230199394865
$ws.Range("U4").Value = 0.9844372815662458
$ws.Range("V4").Value = 0.9851848908345577
$ws.Range("W4").Value = 0.9904860875007423

# New row 5: Rizzie Spiral
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.039738524628567
$ws.Range("D5").Value = 0.8082073620642399
$ws.Range("E5").Value = 1.429159006619342
$ws.Range("F5").Value = 0.8392586835029057
$ws.Range("G5").Value = 1.039738524628567
$ws.Range("H5").Value = 1.108700854873764
$ws.Range("I5").Value = 0.9519489806454151
$ws.Range("J5").Value = 1.429159006619342
$ws.Range("K5").Value = 1.429159006619342
$ws.Range("L5").Value = 0.9989272406430745
$ws.Range("M5").Value = 0.9537681537553273
$ws.Range("N5").Value = 1.429159006619342
$ws.Range("O5").Value = 0.8082073620642399
$ws.Range("P5").Value = 0.9239729433464035
$ws.Range("Q5").Value = 0.8809877579097836
$ws.Range("R5").Value = 1.092368297770716
$ws.Range("S5").Value = 0.9339046801493781
$ws.Range("T5").Value = 1.092368297770716
$ws.Range("U5").Value = 1.057718261766869
$ws.Range("V5").Value = 1.132006410737364
$ws.Range("W5").Value = 1.016213600841579

# Rename "Thomas Hex" to "Matthies Hex" (now located at row 11 after the shift)
$ws.Range("B11").Value = "Matthies Hex"
